# Auto-generated edit script: updates LeveProfit/Price columns (H:N) per the commit diff.
# The workbook stores these as plain numeric literals (no formulas), so we just
# write the new values directly, and clear the handful of cells the diff removes.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 224.375
$ws.Range("I4").Value = 185
$ws.Range("K4").Value = 185
$ws.Range("M4").Value = -71
$ws.Range("H19").Value = 984.7059
$ws.Range("I19").Value = 1047.75
$ws.Range("J19").Value = 833.4
$ws.Range("K19").Value = 1047.75
$ws.Range("L19").Value = 833.4
$ws.Range("M19").Value = -872.75
$ws.Range("N19").Value = -1183.4
$ws.Range("H38").Value = 1069
$ws.Range("I38").Value = 984.4545000000001
$ws.Range("K38").Value = 2953.3635
$ws.Range("M38").Value = -2581.3635
$ws.Range("H41").Value = 1024.7142
$ws.Range("I41").Value = 1194.909
$ws.Range("K41").Value = 1194.909
$ws.Range("M41").Value = -754.9090000000001
$ws.Range("H54").Value = 45000
$ws.Range("I54").Value = 45000
$ws.Range("K54").Value = 45000
$ws.Range("M54").Value = -44514
$ws.Range("H61").Value = 681.0833
$ws.Range("I61").Value = 619.8
$ws.Range("J61").Value = 987.5
$ws.Range("K61").Value = 1859.4
$ws.Range("L61").Value = 2962.5
$ws.Range("M61").Value = -1687.4
$ws.Range("N61").Value = -3306.5
$ws.Range("H64").Value = 3994.6667
$ws.Range("I64").Value = 3992
$ws.Range("J64").Value = 3996
$ws.Range("K64").Value = 3992
$ws.Range("L64").Value = 3996
$ws.Range("M64").Value = -3744
$ws.Range("N64").Value = -4492
$ws.Range("H67").Value = 3994.6667
$ws.Range("I67").Value = 3992
$ws.Range("J67").Value = 3996
$ws.Range("K67").Value = 3992
$ws.Range("L67").Value = 3996
$ws.Range("M67").Value = -3134
$ws.Range("N67").Value = -5712
$ws.Range("H107").Value = 493.8889
$ws.Range("I107").Value = 480.625
$ws.Range("K107").Value = 480.625
$ws.Range("M107").Value = 1439.375
$ws.Range("H113").Value = 4309.5
$ws.Range("I113").Value = 4217.091
$ws.Range("J113").Value = 4454.7144
$ws.Range("K113").Value = 4217.091
$ws.Range("L113").Value = 4454.7144
$ws.Range("M113").Value = -963.0910000000003
$ws.Range("N113").Value = -10962.7144
$ws.Range("H118").Value = 1116.3334
$ws.Range("I118").Value = 339.6
$ws.Range("J118").Value = 5000
$ws.Range("K118").Value = 1018.8
$ws.Range("L118").Value = 15000
$ws.Range("M118").Value = 638.1999999999999
$ws.Range("N118").Value = -18314
$ws.Range("H138").Value = 2562.2068
$ws.Range("J138").Value = 2646.1177
$ws.Range("L138").Value = 7938.353099999999
$ws.Range("N138").Value = -18218.3531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1224.1818
$ws.Range("J2").Value = 1699.1111
$ws.Range("L2").Value = 1699.1111
$ws.Range("N2").Value = -1925.1111
$ws.Range("H21").Value = 2052.3333
$ws.Range("I21").Value = 1805
$ws.Range("J21").Value = 2299.6667
$ws.Range("K21").Value = 1805
$ws.Range("L21").Value = 2299.6667
$ws.Range("M21").Value = -1431
$ws.Range("N21").Value = -3047.6667
$ws.Range("H29").Value = 3597.5
$ws.Range("I29").Value = 3597.5
$ws.Range("K29").Value = 3597.5
$ws.Range("M29").Value = -3289.5
$ws.Range("H45").Value = 3859.2942
$ws.Range("I45").Value = 3838.077
$ws.Range("J45").Value = 3928.25
$ws.Range("K45").Value = 3838.077
$ws.Range("L45").Value = 3928.25
$ws.Range("M45").Value = -3461.077
$ws.Range("N45").Value = -4682.25
$ws.Range("H61").Value = 1811.8
$ws.Range("I61").Value = 1811.8
$ws.Range("K61").Value = 1811.8
$ws.Range("M61").Value = -1599.8
$ws.Range("H94").Value = 535000
$ws.Range("J94").Value = 535000
$ws.Range("L94").Value = 535000
$ws.Range("N94").Value = -536802
$ws.Range("H116").Value = 1224.1818
$ws.Range("J116").Value = 1699.1111
$ws.Range("L116").Value = 1699.1111
$ws.Range("N116").Value = -6287.1111
$ws.Range("H122").Value = 850.6842
$ws.Range("I122").Value = 580.86664
$ws.Range("K122").Value = 1742.59992
$ws.Range("M122").Value = 707.4000800000001
$ws.Range("H132").Value = 8890
$ws.Range("J132").Value = 9340
$ws.Range("L132").Value = 28020
$ws.Range("N132").Value = -33080
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 1811.8
$ws.Range("I136").Value = 1811.8
$ws.Range("K136").Value = 5435.4
$ws.Range("M136").Value = -2885.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1224.1818
$ws.Range("J3").Value = 1699.1111
$ws.Range("L3").Value = 1699.1111
$ws.Range("N3").Value = -1927.1111
$ws.Range("H86").Value = 12275.066
$ws.Range("I86").Value = 13972.6
$ws.Range("K86").Value = 13972.6
$ws.Range("M86").Value = -12849.6
$ws.Range("H89").Value = 12275.066
$ws.Range("I89").Value = 13972.6
$ws.Range("K89").Value = 69863
$ws.Range("M89").Value = -64247
$ws.Range("H99").Value = 1111.9333
$ws.Range("I99").Value = 831.5833
$ws.Range("K99").Value = 831.5833
$ws.Range("M99").Value = 666.4167
$ws.Range("H105").Value = 3012.9092
$ws.Range("I105").Value = 2994.2
$ws.Range("K105").Value = 2994.2
$ws.Range("M105").Value = -1247.2
$ws.Range("H134").Value = 1232.1428
$ws.Range("I134").Value = 1232.1428
$ws.Range("K134").Value = 3696.4284
$ws.Range("M134").Value = -1161.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1499
$ws.Range("I22").Value = 223
$ws.Range("J22").Value = 2775
$ws.Range("K22").Value = 223
$ws.Range("L22").Value = 2775
$ws.Range("M22").Value = 127
$ws.Range("N22").Value = -3475
$ws.Range("H51").Value = 32499.75
$ws.Range("I51").Value = 32499.75
$ws.Range("K51").Value = 32499.75
$ws.Range("M51").Value = -31763.75
$ws.Range("H61").Value = 32499.75
$ws.Range("I61").Value = 32499.75
$ws.Range("K61").Value = 32499.75
$ws.Range("M61").Value = -32151.75
$ws.Range("H99").Value = 1917.2142
$ws.Range("I99").Value = 1989.5714
$ws.Range("J99").Value = 1844.8572
$ws.Range("K99").Value = 1989.5714
$ws.Range("L99").Value = 1844.8572
$ws.Range("M99").Value = -491.5714
$ws.Range("N99").Value = -4840.8572
$ws.Range("H122").Value = 1505
$ws.Range("I122").Value = 1161.75
$ws.Range("J122").Value = 1779.6
$ws.Range("K122").Value = 3485.25
$ws.Range("L122").Value = 5338.799999999999
$ws.Range("M122").Value = -1035.25
$ws.Range("N122").Value = -10238.8
$ws.Range("H126").Value = 1917.2142
$ws.Range("I126").Value = 1989.5714
$ws.Range("J126").Value = 1844.8572
$ws.Range("K126").Value = 5968.7142
$ws.Range("L126").Value = 5534.571599999999
$ws.Range("M126").Value = -3498.7142
$ws.Range("N126").Value = -10474.5716
$ws.Range("H132").Value = 5918.647
$ws.Range("J132").Value = 10000.333
$ws.Range("L132").Value = 30000.999
$ws.Range("N132").Value = -35060.999
$ws.Range("H134").Value = 1554.1765
$ws.Range("I134").Value = 1554.1765
$ws.Range("K134").Value = 4662.529500000001
$ws.Range("M134").Value = -2127.529500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 419.33334
$ws.Range("I17").Value = 310.5
$ws.Range("J17").Value = 506.4
$ws.Range("K17").Value = 931.5
$ws.Range("L17").Value = 1519.2
$ws.Range("M17").Value = -762.5
$ws.Range("N17").Value = -1857.2
$ws.Range("H52").Value = 1933.3334
$ws.Range("J52").Value = 1933.3334
$ws.Range("L52").Value = 5800.0002
$ws.Range("N52").Value = -6332.0002
$ws.Range("H94").Value = 9024.75
$ws.Range("I94").Value = 724.5
$ws.Range("K94").Value = 2173.5
$ws.Range("M94").Value = -1497.5
$ws.Range("H118").Value = 250
$ws.Range("I118").Value = 250
$ws.Range("K118").Value = 750
$ws.Range("M118").Value = 493

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 1500
$ws.Range("I20").Value = 1500
$ws.Range("J20").Value = 1500
$ws.Range("K20").Value = 1500
$ws.Range("L20").Value = 1500
$ws.Range("M20").Value = -1255
$ws.Range("N20").Value = -1990
$ws.Range("H24").Value = 28142.857
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H25").Value = 6168
$ws.Range("I25").Value = 4833
$ws.Range("J25").Value = 7503
$ws.Range("K25").Value = 4833
$ws.Range("L25").Value = 7503
$ws.Range("M25").Value = -4304
$ws.Range("N25").Value = -8561
$ws.Range("H97").Value = 539
$ws.Range("I97").Value = 316.8
$ws.Range("J97").Value = 1650
$ws.Range("K97").Value = 316.8
$ws.Range("L97").Value = 1650
$ws.Range("M97").Value = 179.2
$ws.Range("N97").Value = -2642
$ws.Range("H102").Value = 1708.3704
$ws.Range("I102").Value = 1370.6957
$ws.Range("K102").Value = 1370.6957
$ws.Range("M102").Value = 251.3043
$ws.Range("H126").Value = 4007
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530
$ws.Range("H132").Value = 4992
$ws.Range("I132").Value = 4992
$ws.Range("K132").Value = 14976
$ws.Range("M132").Value = -12446

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1390.7142
$ws.Range("I22").Value = 1393
$ws.Range("K22").Value = 1393
$ws.Range("M22").Value = -1098
$ws.Range("H26").Value = 10000
$ws.Range("J26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("N26").Value = -10590
$ws.Range("H27").Value = 1390.7142
$ws.Range("I27").Value = 1393
$ws.Range("K27").Value = 1393
$ws.Range("M27").Value = -1286
$ws.Range("H40").Value = 3254.8462
$ws.Range("I40").Value = 2946.6365
$ws.Range("K40").Value = 2946.6365
$ws.Range("M40").Value = -2810.6365
$ws.Range("H55").Value = 1427
$ws.Range("I55").Value = 422.5
$ws.Range("K55").Value = 422.5
$ws.Range("M55").Value = -249.5
$ws.Range("H100").Value = 1892.4
$ws.Range("I100").Value = 1892.4
$ws.Range("K100").Value = 1892.4
$ws.Range("M100").Value = -1351.4
$ws.Range("H106").Value = 17180.75
$ws.Range("J106").Value = 17180.75
$ws.Range("L106").Value = 17180.75
$ws.Range("N106").Value = -19704.75
$ws.Range("H132").Value = 2561.158
$ws.Range("I132").Value = 2333.0322
$ws.Range("K132").Value = 6999.096600000001
$ws.Range("M132").Value = -4469.096600000001
$ws.Range("H136").Value = 2712.4285
$ws.Range("J136").Value = 3499
$ws.Range("L136").Value = 10497
$ws.Range("N136").Value = -15597

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 6377.5
$ws.Range("I23").Value = 3505
$ws.Range("J23").Value = 9250
$ws.Range("K23").Value = 3505
$ws.Range("L23").Value = 9250
$ws.Range("M23").Value = -3276
$ws.Range("N23").Value = -9708
$ws.Range("H70").Value = 1000
$ws.Range("I70").Value = 1000
$ws.Range("K70").Value = 1000
$ws.Range("M70").Value = -685
$ws.Range("H73").Value = 1000
$ws.Range("I73").Value = 1000
$ws.Range("K73").Value = 1000
$ws.Range("M73").Value = 92
$ws.Range("H81").Value = 2140.2727
$ws.Range("I81").Value = 1713.1428
$ws.Range("J81").Value = 2887.75
$ws.Range("K81").Value = 3426.2856
$ws.Range("L81").Value = 5775.5
$ws.Range("M81").Value = -2365.2856
$ws.Range("N81").Value = -7897.5
$ws.Range("H84").Value = 2140.2727
$ws.Range("I84").Value = 1713.1428
$ws.Range("J84").Value = 2887.75
$ws.Range("K84").Value = 17131.428
$ws.Range("L84").Value = 28877.5
$ws.Range("M84").Value = -11827.428
$ws.Range("N84").Value = -39485.5
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H132").Value = 5498
$ws.Range("I132").Value = 6997.3335
$ws.Range("K132").Value = 20992.0005
$ws.Range("M132").Value = -18462.0005
$ws.Range("H133").Value = 145000
$ws.Range("J133").Value = 145000
$ws.Range("L133").Value = 145000
$ws.Range("N133").Value = -155120
